# Refresh the cryptos table with the latest Coinranking snapshot values
# (commit: "Updated cryptos list on Thu May 23 21:19:17 UTC 2024 with GitHub Actions").
# Row 18/19, 24/25/26 and 46/47 also swap rank order as coins moved in the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores plain-text numbers (e.g. "1.00", "598.80",
# "67.972.52"). Force Text format on the cells we are about to rewrite with a
# numeric-looking string so Excel keeps the literal text (trailing zeros /
# dotted thousands separators) instead of silently coercing it to a Number.
$ws.Range("D4:D6").NumberFormat = "@"
$ws.Range("D8:D14").NumberFormat = "@"
$ws.Range("D18:D28").NumberFormat = "@"
$ws.Range("D30:D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49:D51").NumberFormat = "@"

$ws.Range("D2").Value = '67.972.52'
$ws.Range("E2").Value = '  -1.66%  '

$ws.Range("D3").Value = '3.836.34'
$ws.Range("E3").Value = '  +2.78%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.19%  '

$ws.Range("D5").Value = '598.80'
$ws.Range("E5").Value = '  -2.19%  '

$ws.Range("D6").Value = '174.80'
$ws.Range("E6").Value = '  -1.21%  '

$ws.Range("D7").Value = '3.840.38'
$ws.Range("E7").Value = '  +2.95%  '

$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("D9").Value = '0.520'
$ws.Range("E9").Value = '  -0.72%  '

$ws.Range("D10").Value = '0.159'
$ws.Range("E10").Value = '  -3.84%  '

$ws.Range("D11").Value = '6.17'
$ws.Range("E11").Value = '  -5.49%  '

$ws.Range("D12").Value = '0.462'
$ws.Range("E12").Value = '  -3.29%  '

$ws.Range("D13").Value = '38.14'
$ws.Range("E13").Value = '  -3.86%  '

$ws.Range("D14").Value = '0.0000246'
$ws.Range("E14").Value = '  -2.23%  '

$ws.Range("D15").Value = '4.475.56'
$ws.Range("E15").Value = '  +2.74%  '

$ws.Range("D16").Value = '3.842.90'
$ws.Range("E16").Value = '  +2.87%  '

$ws.Range("D17").Value = '68.101.04'
$ws.Range("E17").Value = '  -1.69%  '

$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").Value = '7.19'
$ws.Range("E18").Value = '  -2.79%  '

$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.115'
$ws.Range("E19").Value = '  -4.59%  '

$ws.Range("D20").Value = '16.29'
$ws.Range("E20").Value = '  +0.38%  '

$ws.Range("D21").Value = '489.37'
$ws.Range("E21").Value = '  -1.62%  '

$ws.Range("D22").Value = '9.20'
$ws.Range("E22").Value = '  +0.61%  '

$ws.Range("D23").Value = '0.733'
$ws.Range("E23").Value = '  +2.33%  '

$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = '0.0000149'
$ws.Range("E24").Value = '  +12.55%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '83.72'
$ws.Range("E25").Value = '  -2.03%  '

$ws.Range("B26").Value = 'Fetch.AI'
$ws.Range("C26").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D26").Value = '2.37'
$ws.Range("E26").Value = '  -5.97%  '

$ws.Range("D27").Value = '12.26'
$ws.Range("E27").Value = '  -4.32%  '

$ws.Range("D28").Value = '10.24'
$ws.Range("E28").Value = '  -5.27%  '

$ws.Range("E29").Value = '  +0.42%  '

$ws.Range("D30").Value = '2.96'
$ws.Range("E30").Value = '  +2.28%  '

$ws.Range("D31").Value = '33.12'
$ws.Range("E31").Value = '  +9.73%  '

$ws.Range("D32").Value = '2.43'
$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").Value = '7.71'
$ws.Range("E33").Value = '  -3.31%  '

$ws.Range("D34").Value = '0.109'
$ws.Range("E34").Value = '  -2.68%  '

$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").Value = '1.01'
$ws.Range("E36").Value = '  -3.12%  '

$ws.Range("D37").Value = '0.137'
$ws.Range("E37").Value = '  -0.44%  '

$ws.Range("D38").Value = '5.76'
$ws.Range("E38").Value = '  -5.20%  '

$ws.Range("D39").Value = '463.81'
$ws.Range("E39").Value = '  +4.70%  '

$ws.Range("D40").Value = '0.330'
$ws.Range("E40").Value = '  -4.10%  '

$ws.Range("D41").Value = '49.10'
$ws.Range("E41").Value = '  -1.08%  '

$ws.Range("D42").Value = '2.01'
$ws.Range("E42").Value = '  -2.05%  '

$ws.Range("D43").Value = '2.87'
$ws.Range("E43").Value = '  -5.87%  '

$ws.Range("D44").Value = '8.32'
$ws.Range("E44").Value = '  -2.11%  '

$ws.Range("D45").Value = '41.50'
$ws.Range("E45").Value = '  -5.74%  '

$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.847.40'
$ws.Range("E46").Value = '  -3.15%  '

$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D47").Value = '140.86'
$ws.Range("E47").Value = '  +2.32%  '

$ws.Range("D49").Value = '0.0352'
$ws.Range("E49").Value = '  -1.46%  '

$ws.Range("D50").Value = '26.18'
$ws.Range("E50").Value = '  -2.63%  '

$ws.Range("D51").Value = '23.76'
$ws.Range("E51").Value = '  +11.24%  '
